$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a set of cells to remain text-typed when assigned numeric-looking strings,
# then reset the style back to Normal so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.40"
Set-TextValue "G2" "20"
Set-TextValue "D3" "23.38"
Set-TextValue "G3" "20"
Set-TextValue "D4" "5.659"
Set-TextValue "G4" "20"
Set-TextValue "D5" "0.05816"
Set-TextValue "G5" "20"
Set-TextValue "D6" "3.404"
Set-TextValue "G6" "20"
Set-TextValue "D7" "6.467"
Set-TextValue "G7" "20"
Set-TextValue "D8" "1.319"
Set-TextValue "G8" "20"
Set-TextValue "D9" "0.7977"
Set-TextValue "G9" "20"
Set-TextValue "D10" "0.1459"
Set-TextValue "G10" "20"
Set-TextValue "D11" "0.07634"
Set-TextValue "G11" "20"
Set-TextValue "D12" "0.03210"
Set-TextValue "G12" "20"
Set-TextValue "D13" "0.02959"
Set-TextValue "G13" "20"
Set-TextValue "G14" "20"
Set-TextValue "D15" "0.001680"
Set-TextValue "G15" "20"
Set-TextValue "D16" "3.292"
Set-TextValue "G16" "20"
Set-TextValue "D17" "0.04745"
Set-TextValue "G17" "20"
Set-TextValue "D18" "0.0005984"
Set-TextValue "G18" "20"
Set-TextValue "D19" "0.006201"
Set-TextValue "G19" "20"
Set-TextValue "D20" "0.005404"
Set-TextValue "G20" "20"
Set-TextValue "D21" "0.001064"
Set-TextValue "G21" "20"
Set-TextValue "G22" "20"
Set-TextValue "D23" "3.697"
Set-TextValue "G23" "20"
Set-TextValue "G24" "20"
Set-TextValue "D25" "0.3324"
Set-TextValue "G25" "20"
Set-TextValue "D26" "0.1239"
Set-TextValue "G26" "20"
Set-TextValue "D27" "0.0009987"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
Set-TextValue "G27" "20"
Set-TextValue "G28" "20"
Set-TextValue "G29" "20"
Set-TextValue "G30" "20"
Set-TextValue "G31" "20"
Set-TextValue "G32" "20"
Set-TextValue "G33" "20"
Set-TextValue "G34" "20"
Set-TextValue "G35" "20"
Set-TextValue "G36" "20"
Set-TextValue "G37" "20"
Set-TextValue "G38" "20"
Set-TextValue "G39" "20"
Set-TextValue "D40" "0.04294"
Set-TextValue "G40" "20"
Set-TextValue "D41" "0.007075"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue "G41" "20"
Set-TextValue "D42" "0.003598"
Set-TextValue "G42" "20"
Set-TextValue "D43" "0.1057"
Set-TextValue "G43" "20"
Set-TextValue "D44" "0.008757"
Set-TextValue "G44" "20"
Set-TextValue "G45" "20"
Set-TextValue "D46" "0.00005422"
Set-TextValue "G46" "20"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "G47" "20"
Set-TextValue "D48" "0.7844"
Set-TextValue "G48" "20"
Set-TextValue "D49" "0.1005"
Set-TextValue "G49" "20"
Set-TextValue "D50" "0.00002099"
Set-TextValue "G50" "20"
Set-TextValue "G51" "20"
